$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Rushing
$ws2 = $wb.Worksheets.Item(2)   # Receiving

# ---------------------------------------------------------------------------
# Sheet 1 "Rushing": log Week 16 stats.
# A new player (C.Henne) appears and is inserted as row 3, pushing the
# existing rows (old 3-11) down to (4-12). We do this by copying values from
# the bottom up (reading with .Value2, which works reliably; writing with
# .Value) rather than using Rows.Insert(), because Insert() on this engine
# spuriously introduces a brand new (unused) cell style.
# ---------------------------------------------------------------------------

for ($r = 11; $r -ge 3; $r--) {
    $dst = $r + 1
    $ws1.Range("A$dst").Value = $ws1.Range("A$r").Value2
    $ws1.Range("B$dst").Value = $ws1.Range("B$r").Value2
    $ws1.Range("C$dst").Value = $ws1.Range("C$r").Value2
    $ws1.Range("D$dst").Value = $ws1.Range("D$r").Value2
    $ws1.Range("E$dst").Value = $ws1.Range("E$r").Value2
    $ws1.Range("F$dst").Value = $ws1.Range("F$r").Value2
}

# New row 3: C.Henne
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "C.Henne"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 1
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 2

# Match formatting of column A in the new row (and the newly-appended row 12,
# whose A-cell was blank before the shift) to the existing header/index style.
foreach ($addr in @("A3", "A12")) {
    $cell = $ws1.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Updated cumulative rushing totals after Week 16 (rows shifted down by one)
$ws1.Range("C4").Value = 73
$ws1.Range("D4").Value = 44
$ws1.Range("E4").Value = 2
$ws1.Range("F4").Value = 14

$ws1.Range("C5").Value = 39
$ws1.Range("D5").Value = 19
$ws1.Range("E5").Value = 7
$ws1.Range("F5").Value = 16

$ws1.Range("C7").Value = 15
$ws1.Range("D7").Value = 9
$ws1.Range("E7").Value = 2
$ws1.Range("F7").Value = 2

$ws1.Range("C8").Value = 0
$ws1.Range("D8").Value = 0
$ws1.Range("E8").Value = 4
$ws1.Range("F8").Value = 3

# ---------------------------------------------------------------------------
# Sheet 2 "Receiving": season simulation results from Week 17 onward.
# ---------------------------------------------------------------------------

# T.Hill (row 7)
$ws2.Range("C7").Value = 111
$ws2.Range("D7").Value = 88
$ws2.Range("E7").Value = 32
$ws2.Range("F7").Value = 13
$ws2.Range("G7").Value = 21
$ws2.Range("H7").Value = 16

# M.Hardman (row 8)
$ws2.Range("C8").Value = 53
$ws2.Range("D8").Value = 41
$ws2.Range("E8").Value = 13
$ws2.Range("F8").Value = 5

# B.Pringle (row 9)
$ws2.Range("C9").Value = 28
$ws2.Range("D9").Value = 19
$ws2.Range("E9").Value = 13
$ws2.Range("F9").Value = 9

# D.Robinson (row 10)
$ws2.Range("C10").Value = 22
$ws2.Range("D10").Value = 15
$ws2.Range("E10").Value = 7
$ws2.Range("F10").Value = 2

# Row 12 becomes J.Gordon (unchanged stats), row 13 becomes T.Kelce (updated
# cumulative totals).
$ws2.Range("B12").Value = "J.Gordon"

$ws2.Range("B13").Value = "T.Kelce"
$ws2.Range("C13").Value = 98
$ws2.Range("D13").Value = 69
$ws2.Range("E13").Value = 23
$ws2.Range("F13").Value = 13
$ws2.Range("G13").Value = 14
$ws2.Range("H13").Value = 10

# Update the selected cell shown on the Receiving tab
$ws2.Activate()
$ws2.Range("N9").Select()
